$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (619-625) appended to the Covid-19 daily tracking table.
# Column A holds the "Data" (date) label, stored as plain text (matches
# all the existing rows, which are shared strings rather than real dates).
$dates = @(
    "2021/12/07",
    "2021/12/08",
    "2021/12/09",
    "2021/12/10",
    "2021/12/11",
    "2021/12/12",
    "2021/12/13"
)

# Columns B..T for rows 619..625, in order.
$data = @(
    @(618, 21593, 611, 10, 0.0282962071,  20911, 27040, 48633, 1, 0, 1, 2, 0, 2, 4, 3, 11, 61, 89),
    @(619, 21593, 611, 10, 0.0282962071,  20911, 27040, 48633, 0, 0, 0, 2, 0, 2, 4, 3, 11, 61, 89),
    @(620, 21593, 611, 9,  0.0282962071,  20912, 27040, 48633, 0, 0, 0, 2, 0, 2, 5, 2, 11, 61, 89),
    @(621, 21593, 611, 5,  0.0282962071,  20916, 27040, 48633, 0, 0, 0, 1, 0, 1, 3, 2, 11, 61, 89),
    @(622, 21593, 611, 5,  0.0282962071,  20916, 27040, 48633, 0, 0, 0, 1, 0, 1, 3, 2, 11, 61, 89),
    @(623, 21593, 611, 5,  0.0282962071,  20916, 27040, 48633, 0, 0, 0, 1, 0, 1, 3, 2, 11, 61, 90),
    @(624, 21594, 611, 5,  0.02829489673, 20917, 27040, 48634, 1, 0, 1, 1, 0, 1, 4, 1, 11, 61, 90)
)

$startRow = 619
for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = $startRow + $i

    # Write column A as literal text (not auto-converted to a date serial).
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $dates[$i]
    $cellA.Style = "Normal"

    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, 2 + $c).Value = $values[$c]
    }
}
